$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.969.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '1.906.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.17%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = '''1.005'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.79%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''315.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.4805'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.3798'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.07354'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.9314'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''20.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''0.07758'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '1.925.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.69%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''5.494'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''6.634'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''91.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.40%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.000008834'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '28.006.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.63%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''14.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''5.165'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '2.154.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = '''1.921'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = 'Monero'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = '''155.45'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '  +5.29%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''116.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.18%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''4.952'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.08938'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.311'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''1.257'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''0.7739'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''4.673'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Value = '''2.627'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.77%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.02053'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''1.111'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.05307'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'MXToken'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '''2.997'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''0.5477'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''7.014'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.1522'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''8.468'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''10.63'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''0.4824'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''108.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.01%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '  -0.90%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''1.648'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.38%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''67.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.06072'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.12%  '
$ws.Range('E51').Style = 'Normal'
